$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "28,57 TL - 28,57 TL"
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F7").Value = "%3"
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

$ws.Range("F14").Value = ""
